# Auto-generated edit script: update market-price / profit cells per commit
# "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 17168.812
$ws.Range("I74").Value = 17713.133
$ws.Range("K74").Value = 17713.133
$ws.Range("M74").Value = -16777.133
$ws.Range("H77").Value = 17168.812
$ws.Range("I77").Value = 17713.133
$ws.Range("K77").Value = 88565.66500000001
$ws.Range("M77").Value = -83885.66500000001
$ws.Range("H106").Value = 3561.875
$ws.Range("I106").Value = 3613.5715
$ws.Range("K106").Value = 3613.5715
$ws.Range("M106").Value = -2982.5715
$ws.Range("H132").Value = 5123.5625
$ws.Range("I132").Value = 5754.148
$ws.Range("K132").Value = 17262.444
$ws.Range("M132").Value = -14732.444
$ws.Range("H138").Value = 2838.0112
$ws.Range("I138").Value = 2294.6667
$ws.Range("J138").Value = 3070.873
$ws.Range("K138").Value = 6884.000100000001
$ws.Range("L138").Value = 9212.619000000001
$ws.Range("M138").Value = -1744.000100000001
$ws.Range("N138").Value = -19492.619
$ws.Range("H141").Value = 1170.25
$ws.Range("I141").Value = 908.8570999999999
$ws.Range("K141").Value = 2726.5713
$ws.Range("M141").Value = 2453.4287

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1952
$ws.Range("J2").Value = 1999.5
$ws.Range("L2").Value = 1999.5
$ws.Range("N2").Value = -2225.5
$ws.Range("H32").Value = 4564.38
$ws.Range("I32").Value = 3898.4468
$ws.Range("J32").Value = 14997.333
$ws.Range("K32").Value = 3898.4468
$ws.Range("L32").Value = 14997.333
$ws.Range("M32").Value = -3611.4468
$ws.Range("N32").Value = -15571.333
$ws.Range("H45").Value = 15032.294
$ws.Range("I45").Value = 31950.357
$ws.Range("K45").Value = 31950.357
$ws.Range("M45").Value = -31573.357
$ws.Range("H102").Value = 4516.5293
$ws.Range("I102").Value = 4548.6787
$ws.Range("K102").Value = 4548.6787
$ws.Range("M102").Value = -2926.6787
$ws.Range("H116").Value = 1952
$ws.Range("J116").Value = 1999.5
$ws.Range("L116").Value = 1999.5
$ws.Range("N116").Value = -6587.5
$ws.Range("H122").Value = 5130.727
$ws.Range("I122").Value = 4156.6665
$ws.Range("K122").Value = 12469.9995
$ws.Range("M122").Value = -10019.9995
$ws.Range("H132").Value = 2807.9546
$ws.Range("I132").Value = 2461.8667
$ws.Range("J132").Value = 3549.5715
$ws.Range("K132").Value = 7385.6001
$ws.Range("L132").Value = 10648.7145
$ws.Range("M132").Value = -4855.6001
$ws.Range("N132").Value = -15708.7145

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1952
$ws.Range("J3").Value = 1999.5
$ws.Range("L3").Value = 1999.5
$ws.Range("N3").Value = -2227.5
$ws.Range("H105").Value = 7431963.5
$ws.Range("I105").Value = 387895.53
$ws.Range("K105").Value = 387895.53
$ws.Range("M105").Value = -386148.53
$ws.Range("H107").Value = 1564.909
$ws.Range("J107").Value = 1966.6666
$ws.Range("L107").Value = 1966.6666
$ws.Range("N107").Value = -5806.6666
$ws.Range("H134").Value = 3609.9412
$ws.Range("I134").Value = 3320.111
$ws.Range("J134").Value = 3936
$ws.Range("K134").Value = 9960.332999999999
$ws.Range("L134").Value = 11808
$ws.Range("M134").Value = -7425.332999999999
$ws.Range("N134").Value = -16878

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 58827260
$ws.Range("I7").Value = 4619.154
$ws.Range("J7").Value = 250000850
$ws.Range("K7").Value = 4619.154
$ws.Range("L7").Value = 250000850
$ws.Range("M7").Value = -4506.154
$ws.Range("N7").Value = -250001076
$ws.Range("H31").Value = 3287.1733
$ws.Range("I31").Value = 2792.5
$ws.Range("J31").Value = 4559.1904
$ws.Range("K31").Value = 2792.5
$ws.Range("L31").Value = 4559.1904
$ws.Range("M31").Value = -2497.5
$ws.Range("N31").Value = -5149.1904
$ws.Range("H34").Value = 3287.1733
$ws.Range("I34").Value = 2792.5
$ws.Range("J34").Value = 4559.1904
$ws.Range("K34").Value = 2792.5
$ws.Range("L34").Value = 4559.1904
$ws.Range("M34").Value = -2590.5
$ws.Range("N34").Value = -4963.1904
$ws.Range("H58").Value = 3913.6667
$ws.Range("J58").Value = 4925.8
$ws.Range("L58").Value = 4925.8
$ws.Range("N58").Value = -5331.8
$ws.Range("H92").Value = 39981.5
$ws.Range("J92").Value = 39981.5
$ws.Range("L92").Value = 39981.5
$ws.Range("N92").Value = -44973.5
$ws.Range("H105").Value = 1972.55
$ws.Range("I105").Value = 1300.9
$ws.Range("J105").Value = 2644.2
$ws.Range("K105").Value = 1300.9
$ws.Range("L105").Value = 2644.2
$ws.Range("M105").Value = 446.0999999999999
$ws.Range("N105").Value = -6138.2
$ws.Range("H107").Value = 992.3570999999999
$ws.Range("I107").Value = 1109.6
$ws.Range("J107").Value = 699.25
$ws.Range("K107").Value = 1109.6
$ws.Range("L107").Value = 699.25
$ws.Range("M107").Value = 810.4000000000001
$ws.Range("N107").Value = -4539.25
$ws.Range("H134").Value = 2784.9487
$ws.Range("I134").Value = 2174.2222
$ws.Range("J134").Value = 4159.0835
$ws.Range("K134").Value = 6522.6666
$ws.Range("L134").Value = 12477.2505
$ws.Range("M134").Value = -3987.6666
$ws.Range("N134").Value = -17547.2505
$ws.Range("H136").Value = 3913.6667
$ws.Range("J136").Value = 4925.8
$ws.Range("L136").Value = 14777.4
$ws.Range("N136").Value = -19877.4
$ws.Range("H141").Value = 382921.16
$ws.Range("J141").Value = 382921.16
$ws.Range("L141").Value = 382921.16
$ws.Range("N141").Value = -393281.16

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 905.5217
$ws.Range("J107").Value = 923.1579
$ws.Range("L107").Value = 2769.4737
$ws.Range("N107").Value = -6609.4737
$ws.Range("H139").Value = 3701.25
$ws.Range("J139").Value = 3999.6924
$ws.Range("L139").Value = 11999.0772
$ws.Range("N139").Value = -22279.0772
$ws.Range("H140").Value = 12294.179
$ws.Range("I140").Value = 7077.375
$ws.Range("K140").Value = 21232.125
$ws.Range("M140").Value = -16052.125
$ws.Range("H141").Value = 33111
$ws.Range("I141").Value = 19333
$ws.Range("K141").Value = 57999
$ws.Range("M141").Value = -52819

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 250002990
$ws.Range("I80").Value = 500001500
$ws.Range("J80").Value = 4497.5
$ws.Range("K80").Value = 500001500
$ws.Range("L80").Value = 4497.5
$ws.Range("M80").Value = -500000502
$ws.Range("N80").Value = -6493.5
$ws.Range("H83").Value = 250002990
$ws.Range("I83").Value = 500001500
$ws.Range("J83").Value = 4497.5
$ws.Range("K83").Value = 2500007500
$ws.Range("L83").Value = 22487.5
$ws.Range("M83").Value = -2500002508
$ws.Range("N83").Value = -32471.5
$ws.Range("H97").Value = 616.2222
$ws.Range("I97").Value = 616.2222
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 616.2222
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -120.2222
$ws.Range("N97").ClearContents()
$ws.Range("H107").Value = 458.2857
$ws.Range("I107").Value = 458.2857
$ws.Range("K107").Value = 458.2857
$ws.Range("M107").Value = 1461.7143
$ws.Range("H132").Value = 4452.7812
$ws.Range("I132").Value = 3812.8696
$ws.Range("J132").Value = 6088.1113
$ws.Range("K132").Value = 11438.6088
$ws.Range("L132").Value = 18264.3339
$ws.Range("M132").Value = -8908.6088
$ws.Range("N132").Value = -23324.3339
$ws.Range("H134").Value = 67924.75
$ws.Range("J134").Value = 67924.75
$ws.Range("L134").Value = 203774.25
$ws.Range("N134").Value = -208844.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H22").Value = 119050250
$ws.Range("I22").Value = 10205875
$ws.Range("K22").Value = 10205875
$ws.Range("M22").Value = -10205580
$ws.Range("H27").Value = 119050250
$ws.Range("I27").Value = 10205875
$ws.Range("K27").Value = 10205875
$ws.Range("M27").Value = -10205768
$ws.Range("H46").Value = 1522.4103
$ws.Range("J46").Value = 1400.6666
$ws.Range("L46").Value = 1400.6666
$ws.Range("N46").Value = -1776.6666
$ws.Range("H61").Value = 5748.75
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 7331.6665
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 7331.6665
$ws.Range("M61").Value = -798
$ws.Range("N61").Value = -7735.6665
$ws.Range("H113").Value = 5748.75
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 7331.6665
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 7331.6665
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -11671.6665
$ws.Range("H122").Value = 6025.143
$ws.Range("I122").Value = 6221.143
$ws.Range("J122").Value = 5829.143
$ws.Range("K122").Value = 18663.429
$ws.Range("L122").Value = 17487.429
$ws.Range("M122").Value = -16213.429
$ws.Range("N122").Value = -22387.429
$ws.Range("H136").Value = 6853.7085
$ws.Range("I136").Value = 6062.4375
$ws.Range("K136").Value = 18187.3125
$ws.Range("M136").Value = -15637.3125

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4043.5454
$ws.Range("I132").Value = 3758.9792
$ws.Range("J132").Value = 5994.857
$ws.Range("K132").Value = 11276.9376
$ws.Range("L132").Value = 17984.571
$ws.Range("M132").Value = -8746.937600000001
$ws.Range("N132").Value = -23044.571
$ws.Range("H140").Value = 134760.33
$ws.Range("J140").Value = 134760.33
$ws.Range("L140").Value = 134760.33
$ws.Range("N140").Value = -145120.33
$ws.Range("H141").Value = 69544.55
$ws.Range("J141").Value = 69544.55
$ws.Range("L141").Value = 69544.55
$ws.Range("N141").Value = -79904.55
